# The document contains two "<id>...</id>" blocks (for p145r_1 and p145r_2)
# that were previously split across three runs:
#   1) "<id>"     -- Courier New, color 7f6000, sz 18
#   2) "p145r_N"  -- color 000000 (plain)
#   3) "</id>"    -- Courier New, color 7f6000, sz 18
#
# The edit collapses each trio into a single run containing the full
# "<id>p145r_N</id>" text, using the formatting of the first ("<id>") run
# (Courier New / 7f6000 / sz 18). A sibling "<id>fig_p145r_1</id>" run
# elsewhere in the document must be left untouched.
#
# Using Find/Replace scoped to the specific paragraph keeps the edit
# targeted: Word merges the matched text into a single run carrying the
# formatting of the first character of the match, which is exactly the
# "<id>" run's formatting.

$d = $word.ActiveDocument

# --- Block 1: "<id>p145r_1</id>" (paragraph 5: " <id>p145r_1</id>") ---
$rng1 = $d.Paragraphs(5).Range
$rng1.Find.Execute("<id>p145r_1</id>", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "<id>p145r_1</id>", 2)

# --- Block 2: "<id>p145r_2</id>" (paragraph 25: "<id>p145r_2</id>") ---
$rng2 = $d.Paragraphs(25).Range
$rng2.Find.Execute("<id>p145r_2</id>", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "<id>p145r_2</id>", 2)
